$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "david" row (row 3) entirely - remaining rows shift up.
$ws.Rows("3").Delete()

# Update D2 (charlie's height) value.
$ws.Range("D2").Value = 30303030

# Update B3 (now helen's row, after the shift) value.
$ws.Range("B3").Value = 10000

# Update selection to match the target workbook state.
$ws.Range("D3").Select()
